$p = $ppt.ActivePresentation

# The deck currently ends with slide 14, the "Fragen?" (Questions?) slide.
# Duplicate it so a new slide 15 is appended with the original "Fragen?"
# content, then repurpose the original slide 14 as the new "LIVE DEMO"
# slide that now precedes it.
$s = $p.Slides.Item(14)
$s.Duplicate() | Out-Null

$titleShape = $s.Shapes.Item(5)
$titleShape.TextFrame.TextRange.Text = "LIVE DEMO"
